# Apply updates to 北京-漫展信息.xlsx per commit "Update gh-pages to output generated at 456a3b4"
# Sheet order in workbook: 1=展览, 2=演出, 3=本地生活, 4=全部类型

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Worksheets.Item(1)) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 5934
$ws1.Range("F17").Value = 175
$ws1.Range("F18").Value = 634
$ws1.Range("F19").Value = 4534
$ws1.Range("F20").Value = 4534
$ws1.Range("I27").Value = "//i1.hdslb.com/bfs/openplatform/202410/rRPjfmy81729134090035.jpeg"
$ws1.Range("F29").Value = 48
$ws1.Range("F38").Value = 1230
$ws1.Range("F39").Value = 1209

# --- Sheet 2: 演出 (Worksheets.Item(2)) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F11").Value = 14
$ws2.Range("F19").Value = 307
$ws2.Range("F21").Value = 497

# --- Sheet 3: 本地生活 (Worksheets.Item(3)) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F3").Value = 714
$ws3.Range("F4").Value = 194
$ws3.Range("F5").Value = 287

# --- Sheet 4: 全部类型 (Worksheets.Item(4)) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 714
$ws4.Range("F5").Value = 194
$ws4.Range("F6").Value = 5934
$ws4.Range("F19").Value = 14
$ws4.Range("F26").Value = 175
$ws4.Range("F27").Value = 634
$ws4.Range("F28").Value = 4534
$ws4.Range("F29").Value = 4534
$ws4.Range("F37").Value = 48
$ws4.Range("F43").Value = 307
$ws4.Range("F45").Value = 497
$ws4.Range("F48").Value = 1230
$ws4.Range("F50").Value = 1209
